$wb = $excel.ActiveWorkbook

# The "NotifyTemplate_Msg" value (column AF, row 2) on every sheet is being
# renamed from " ET_OOFS_Notify " to "ET_ETRS_Notify" to match the test
# machine's notification template.
$newValue = "ET_ETRS_Notify"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Activate()
    $ws.Range("AF2").Value = $newValue
    $ws.Range("AF2").Select()
}

# Leave sheet 1 active again with its own selection, matching the final
# recorded UI state (AA2 selected / scrolled to show column Y onward).
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("AA2").Select()
